# Applies the "Added all results for dataset with 1 line (0.90,0.92,0.95,0.97,0.99)"
# edit to slide 2 of the flow-chart deck:
#   - shrink/raise the two "Generate Noisy images (SaltPepperNoise.py)" process
#     boxes (ids 10 & 28) and their attached connectors (ids 4, 12 & 30) so the
#     boxes sit directly under "Start"/the picture instead of spanning two lines
#   - drop the "(SaltPepperNoise.py)" second line from box id 10
#   - retitle box id 28 to "Run RANSAC"
#
# NB: PowerPoint's Shape.Top/Height (etc.) are expressed in points while the
# underlying OOXML stores EMU (1 pt = 12700 EMU); the literals below carry
# extra decimal digits so that, after the host's internal float rounding,
# they land back on the exact target EMU values from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

# --- Connector 4 (Start -> "Generate Noisy images" box #10) -------------
$shape4 = Get-ShapeById $s 4
$shape4.VerticalFlip = 0
$shape4.Top = 77.28445053100586
$shape4.Height = 0

# --- Flowchart Process 10 (id 10): "Generate Noisy images" box ----------
$shape10 = Get-ShapeById $s 10
$shape10.Top = 52.843976974487305
$shape10.Height = 48.88090567371038
$shape10.TextFrame.TextRange.Text = "Generate Noisy images "

# --- Connector 12 ("Generate Noisy images" box #10 -> box #11) ----------
$shape12 = Get-ShapeById $s 12
$shape12.Top = 77.28445053100586
$shape12.Height = 0

# --- Flowchart Process 27 (id 28): retitled "Run RANSAC" box ------------
$shape28 = Get-ShapeById $s 28
$shape28.Top = 187.0459442138672
$shape28.Height = 48.88059043884277
$shape28.TextFrame.TextRange.Text = "Run RANSAC"

# --- Connector 30 (Picture 19 -> "Run RANSAC" box #28) -------------------
$shape30 = Get-ShapeById $s 30
$shape30.Height = 55.61870193481477
